# Generate Report for Handback
# Applies the "handback" update to localization-status.xlsx:
#  - Overview/zh-cn/de-de "Status"-ish cells move from "In Translation"
#    to "Handed back: in sync with en-US"
#  - zh-cn & de-de "Latest Handback DateTime" columns get a real timestamp
#  - zh-cn & de-de "Latest Target File" (I) / "Latest Handback File" (J)
#    columns get populated (I as a hyperlink to the source .md)
#  - a handful of columns are widened so the new content is readable

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status text: "In Translation" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: populate "Latest Target File" (I) / "Latest Handback
#    File" (J) and bump "Latest Handback DateTime" (K)
# ---------------------------------------------------------------------
$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/141aa2ba435b98265709e201d587f52c3996edf4/e2e/046b879a-113a-4cb9-a390-9ca6f0a97c9f.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/141aa2ba435b98265709e201d587f52c3996edf4/e2e/25ed901c-aa0d-4c84-8efb-cc5e45349758.md"
$mdName1 = "046b879a-113a-4cb9-a390-9ca6f0a97c9f.md"
$mdName2 = "25ed901c-aa0d-4c84-8efb-cc5e45349758.md"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl1, "", "", $mdName1)
$wsZhCn.Range("I2").Font.Underline = 2
$wsZhCn.Range("I2").Font.Color = 15570276
$wsZhCn.Range("J2").Value = "046b879a-113a-4cb9-a390-9ca6f0a97c9f.643d6614d321e589c589f3e524924b29c35fe44a.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-16 16:21:33"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $mdUrl2, "", "", $mdName2)
$wsZhCn.Range("I3").Font.Underline = 2
$wsZhCn.Range("I3").Font.Color = 15570276
$wsZhCn.Range("J3").Value = "25ed901c-aa0d-4c84-8efb-cc5e45349758.746ad7dd56c3885b3f7f2588f99f2f6ed19f8bc5.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-16 16:21:33"

# ---------------------------------------------------------------------
# 3. de-de sheet: same shape, different handback timestamp
# ---------------------------------------------------------------------
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl1, "", "", $mdName1)
$wsDeDe.Range("I2").Font.Underline = 2
$wsDeDe.Range("I2").Font.Color = 15570276
$wsDeDe.Range("J2").Value = "046b879a-113a-4cb9-a390-9ca6f0a97c9f.643d6614d321e589c589f3e524924b29c35fe44a.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-16 16:21:39"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $mdUrl2, "", "", $mdName2)
$wsDeDe.Range("I3").Font.Underline = 2
$wsDeDe.Range("I3").Font.Color = 15570276
$wsDeDe.Range("J3").Value = "25ed901c-aa0d-4c84-8efb-cc5e45349758.746ad7dd56c3885b3f7f2588f99f2f6ed19f8bc5.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-16 16:21:39"

# ---------------------------------------------------------------------
# 4. Widen columns so the newly-populated cells are readable
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664
